$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ЛС")

# Fill in the remaining details for the existing row 5 (person #4)
$ws.Range("J5").Value = "Манси"
$ws.Range("K5").Value = "М"

# Add new row 6 for person #5
$ws.Range("A6").Value = 5
$ws.Range("F6").Value = "Токаев Аксай Берметович"
$ws.Range("B6").Value = "ЛК-884912"
$ws.Range("C6").Value = "рядовой"

# Widen column F (ФИО) so the new, longer name fits, as Excel auto-fit would do
$ws.Columns.Item(6).ColumnWidth = 38.47

# Select B6 as the active cell, mirroring the resulting selection state
$ws.Range("B6").Select()
